# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.841.41"
$ws.Range("E2").Value = "  +9.44%  "
$ws.Range("D3").Value = "2.687.38"
$ws.Range("E3").Value = "  +10.92%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'188.01"
$ws.Range("E5").Value = "  +13.83%  "
$ws.Range("D6").Value = "'588.26"
$ws.Range("E6").Value = "  +4.58%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D8").Value = "'0.540"
$ws.Range("E8").Value = "  +5.74%  "
$ws.Range("D9").Value = "'0.197"
$ws.Range("E9").Value = "  +17.85%  "
$ws.Range("D10").Value = "2.684.90"
$ws.Range("E10").Value = "  +10.93%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  +7.93%  "
$ws.Range("D13").Value = "'4.73"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "75.593.55"
$ws.Range("E14").Value = "  +9.28%  "
$ws.Range("D15").Value = "3.170.07"
$ws.Range("E15").Value = "  +10.31%  "
$ws.Range("D16").Value = "'0.0000189"
$ws.Range("E16").Value = "  +6.67%  "
$ws.Range("D17").Value = "'26.59"
$ws.Range("E17").Value = "  +11.60%  "
$ws.Range("D18").Value = "2.694.91"
$ws.Range("E18").Value = "  +11.36%  "
$ws.Range("D19").Value = "'9.27"
$ws.Range("E19").Value = "  +30.66%  "
$ws.Range("E20").Value = "  +11.49%  "
$ws.Range("D21").Value = "'373.19"
$ws.Range("E21").Value = "  +9.88%  "
$ws.Range("D22").Value = "'2.29"
$ws.Range("E22").Value = "  +16.73%  "
$ws.Range("D23").Value = "'4.10"
$ws.Range("E23").Value = "  +5.84%  "
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'70.09"
$ws.Range("E26").Value = "  +6.87%  "
$ws.Range("D27").Value = "'4.19"
$ws.Range("E27").Value = "  +10.39%  "
$ws.Range("D28").Value = "'9.39"
$ws.Range("E28").Value = "  +11.28%  "
$ws.Range("D29").Value = "2.828.66"
$ws.Range("E29").Value = "  +10.78%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").Value = "0.0₃0950"
$ws.Range("E31").Value = "  +12.68%  "
$ws.Range("D32").Value = "'1.42"
$ws.Range("E32").Value = "  +15.64%  "
$ws.Range("D33").Value = "'518.12"
$ws.Range("E33").Value = "  +14.89%  "
$ws.Range("D34").Value = "'7.75"
$ws.Range("E34").Value = "  +5.37%  "
$ws.Range("D35").Value = "'1.76"
$ws.Range("E35").Value = "  +9.23%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'163.22"
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("D38").Value = "'0.119"
$ws.Range("E38").Value = "  +8.17%  "
$ws.Range("D39").Value = "'19.25"
$ws.Range("E39").Value = "  +6.18%  "
$ws.Range("D40").Value = "'19.38"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D42").Value = "'5.02"
$ws.Range("E42").Value = "  +14.81%  "
$ws.Range("D43").Value = "'169.72"
$ws.Range("E43").Value = "  +26.48%  "
$ws.Range("D44").Value = "'1.71"
$ws.Range("E44").Value = "  +13.12%  "
$ws.Range("E45").Value = "  +8.86%  "
$ws.Range("E46").Value = "  +10.58%  "
$ws.Range("D47").Value = "'2.38"
$ws.Range("E47").Value = "  +14.62%  "
$ws.Range("D48").Value = "'39.12"
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("D49").Value = "'0.0845"
$ws.Range("E49").Value = "  +16.55%  "
$ws.Range("E50").Value = "  +8.35%  "
$ws.Range("E51").Value = "  +10.25%  "
